$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape

$ws.Range("D2").Value = "'26.315.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").Value = "'1.610.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'213.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").Value = "'18.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.68%  "

$ws.Range("D12").Value = "'1.833.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "'1.609.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").Value = "'26.270.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "'62.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.28%  "

$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").Value = "'202.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("E21").Value = "  +1.35%  "

$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").Value = "'6.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("E24").Value = "  +1.78%  "

$ws.Range("D25").Value = "'143.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "

$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").Value = "'15.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("D29").Value = "'6.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("E30").Value = "  +5.51%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  +2.98%  "

$ws.Range("D33").Value = "'2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("E34").Value = "  +1.22%  "

$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("D36").Value = "'1.163.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "

$ws.Range("D37").Value = "'0.0168"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.56%  "

$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.791"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.89%  "

$ws.Range("D41").Value = "'0.497"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.02%  "

$ws.Range("D42").Value = "'5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "'1.744.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("D45").Value = "'92.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("E46").Value = "  +14.41%  "

$ws.Range("E47").Value = "  +1.33%  "

$ws.Range("D48").Value = "'53.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("E51").Value = "  -0.33%  "
